$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-10 (row 23)
$ws.Range("B23").Value = 6303
$ws.Range("D23").Value = 5867677
$ws.Range("E23").Value = 930.9339996826908
$ws.Range("F23").Value = 8.150308853809207
$ws.Range("H23").Value = 25.72929900343117
